# new updates for home page
# Rewrites the "Actual" (column L) values in the TC_MiniCart sheet from the old
# generic "<verb> succesfully" / "Text Entered successfully: X" style messages
# into the new "<Action>: <Input1-or-null>" style, and adjusts the row heights
# that Excel recalculated (wrap-text autofit) as a result of the shorter/longer
# text.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("TC_MiniCart")

# row -> new "Actual" text (column L), built from each row's Action (col D) and
# Input1 (col G), following the "{Action}: {Input1 or null}" convention.
$actuals = @{
    2  = "VerifyTitle: City Market Norwalk - Online Grocery Supermarket with Home Delivery"
    3  = "VerifyElement: null"
    4  = "Click: null"
    5  = "VerifyElement: null"
    6  = "SetText: Randomemailid"
    7  = "SetText: 123456"
    8  = "Click: null"
    9  = "VerifyText: Akash sangal"
    10 = "VerifyTitle: City Market Norwalk - Online Grocery Supermarket with Home Delivery"
    11 = "VerifyElement: null"
    12 = "MoveToProductList: Quick & Easy Food Solutions"
    13 = "Click: null"
    14 = "VerifyElement: null"
    15 = "VerifyText: You have no items in your shopping cart."
    16 = "Click: null"
    17 = "VerifyNoElement: null"
    18 = "MoveToProductList: Quick & Easy Food Solutions"
    19 = "MoveAndAddProduct: null"
    20 = "MoveAndAddProduct: null"
    21 = "MoveAndAddProduct: null"
    22 = "Click: null"
    23 = "VerifyMiniCart: You have no items in your shopping cart."
    24 = "Click: null"
    25 = "VerifyNoElement: null"
    26 = "Click: null"
    27 = "VerifyElement: null"
    28 = "Click: null"
    29 = "VerifyNoElement: null"
    30 = "Click: null"
    31 = "AddProductfromMiniCart: null"
    32 = "VerifyMiniCart: You have no items in your shopping cart."
    33 = "RemoveProductfromMiniCart: null"
    34 = "RemoveProductfromMiniCart: null"
    35 = "VerifyMiniCart: You have no items in your shopping cart."
    36 = "DeleteProductfromMiniCart: null"
    37 = "VerifyMiniCart: You have no items in your shopping cart."
    38 = "DeleteProductfromMiniCart: all"
    39 = "Click: null"
    40 = "MoveToProductList: Quick & Easy Food Solutions"
    41 = "MoveAndAddProduct: null"
    42 = "Click: null"
    43 = "VerifyMiniCartMsg: You have no items in your shopping cart."
    44 = "DeleteProductfromMiniCart: all"
    45 = "Click: null"
    46 = "Click: null"
    47 = "Click: null"
    48 = "Wait: 6000"
    49 = "VerifyElement: null"
}

foreach ($row in $actuals.Keys) {
    $ws.Range("L$row").Value = $actuals[$row]
}

# Row heights that Excel's wrap-text autofit recalculated once the "Actual"
# column text changed length.
$rowHeights = @{
    2  = 135
    6  = 75
    7  = 45
    10 = 135
    15 = 90
    16 = 30
    23 = 105
    24 = 30
    28 = 30
    31 = 45
    39 = 30
    45 = 30
    46 = 30
    47 = 30
    48 = 15
}

foreach ($row in $rowHeights.Keys) {
    $ws.Rows.Item($row).RowHeight = $rowHeights[$row]
}
